$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (interest count) in column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 145
$ws1.Range("F4").Value = 112
$ws1.Range("F5").Value = 1242
$ws1.Range("F6").Value = 17532
$ws1.Range("F7").Value = 329
$ws1.Range("F10").Value = 6651
$ws1.Range("F11").Value = 676
$ws1.Range("F14").Value = 53
$ws1.Range("F15").Value = 142
$ws1.Range("F17").Value = 163
$ws1.Range("F22").Value = 29
$ws1.Range("F24").Value = 949
$ws1.Range("F26").Value = 5119
$ws1.Range("F27").Value = 528
$ws1.Range("F28").Value = 59
$ws1.Range("F29").Value = 11805
$ws1.Range("F32").Value = 188
$ws1.Range("F33").Value = 256
$ws1.Range("F35").Value = 285

# Sheet "全部类型" (All types) - same rows (shifted by +1 from row 28 onward
# because this sheet also includes a "演出" row not present in "展览")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 145
$ws4.Range("F4").Value = 112
$ws4.Range("F5").Value = 1242
$ws4.Range("F6").Value = 17532
$ws4.Range("F7").Value = 329
$ws4.Range("F10").Value = 6651
$ws4.Range("F11").Value = 676
$ws4.Range("F14").Value = 53
$ws4.Range("F15").Value = 142
$ws4.Range("F17").Value = 163
$ws4.Range("F22").Value = 29
$ws4.Range("F24").Value = 949
$ws4.Range("F26").Value = 5119
$ws4.Range("F27").Value = 528
$ws4.Range("F29").Value = 59
$ws4.Range("F30").Value = 11805
$ws4.Range("F33").Value = 188
$ws4.Range("F34").Value = 256
$ws4.Range("F36").Value = 285
